# Region VIII_LMS.xlsx edits
#  - widen column U (21st column)
#  - fix header casing B1: REGION -> Region
#  - apply the date/time number format to S:W cells that hold real date
#    serial numbers (so they render as "YYYY-MM-DD HH:MM:SS" instead of
#    plain numbers)
#  - cells in S:W that were placeholder 0 become literal "00:00:00" text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width: min=21 max=21 (column U) 16 -> 21 ---
$ws.Columns.Item(21).ColumnWidth = 20.166666666666668

# --- header rename ---
$ws.Range("B1").Value = "Region"

# --- apply date/time format (style s=6) to non-zero S:W date cells ---
$ws.Range("S2:W3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S4:U8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S16:W21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S33:W36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S38:W44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S46:W54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S57:W83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S84:U84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("S85:W89").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- replace placeholder 0 values in S:W with literal "00:00:00" text ---
$ws.Range("V4:W12").Value = "00:00:00"
$ws.Range("S22:W23").Value = "00:00:00"
$ws.Range("S25:W27").Value = "00:00:00"
$ws.Range("S37:W37").Value = "00:00:00"
$ws.Range("S55:W56").Value = "00:00:00"
$ws.Range("V84:W84").Value = "00:00:00"
